$d = $word.ActiveDocument

$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute("<id>p093v_1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p093v_1</id>", 2)

$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute("<id>p093v_2</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p093v_2</id>", 2)

$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute("<id>p093v_3</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p093v_3</id>", 2)
